$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Home win": 3 data rows -> 5 data rows (rows 2-4 -> rows 2-6)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Home win")

$ws1.Range("A2").Value = "19-12-2024 19:00"
$ws1.Range("B2").Value = "NETHERLANDS"
$ws1.Range("C2").Value = "KNVB BEKER"
$ws1.Range("D2").Value = "Eemdijk - Barendrecht"
$ws1.Range("E2").Value = 70
$ws1.Range("F2").Value = 2.8

$ws1.Range("A3").Value = "19-12-2024 18:00"
$ws1.Range("B3").Value = "SPAIN"
$ws1.Range("C3").Value = "SEGUNDA DIVISIÓN"
$ws1.Range("D3").Value = "Cadiz - Burgos"
$ws1.Range("E3").Value = 70
$ws1.Range("F3").Value = 1.91

$ws1.Range("A4").Value = "20-12-2024 19:45"
$ws1.Range("B4").Value = "ENGLAND"
$ws1.Range("C4").Value = "LEAGUE ONE"
$ws1.Range("D4").Value = "Stockport County - Peterborough"
$ws1.Range("E4").Value = 70
$ws1.Range("F4").Value = 1.76

$ws1.Range("A5").Value = "20-12-2024 13:00"
$ws1.Range("B5").Value = "ISRAEL"
$ws1.Range("C5").Value = "LIGA LEUMIT"
$ws1.Range("D5").Value = "Hapoel Kfar Shalem - Hapoel Ramat HaSharon"
$ws1.Range("E5").Value = 80
$ws1.Range("F5").Value = 1.91

$ws1.Range("A6").Value = "20-12-2024 13:00"
$ws1.Range("B6").Value = "WORLD"
$ws1.Range("C6").Value = "AFF CHAMPIONSHIP"
$ws1.Range("D6").Value = "Malaysia - Singapore"
$ws1.Range("E6").Value = 80
$ws1.Range("F6").Value = 1.7

# ---------------------------------------------------------------
# Sheet "Draw": 3 data rows -> 1 data row (rows 2-4 -> row 2)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Draw")

$ws2.Range("A2").Value = "19-12-2024 18:00"
$ws2.Range("B2").Value = "ROMANIA"
$ws2.Range("C2").Value = "CUPA ROMÂNIEI"
$ws2.Range("D2").Value = "CS Afumati - Arges Pitesti"
$ws2.Range("E2").Value = 66.7
$ws2.Range("F2").Value = 3.1

$ws2.Rows.Item(4).Delete()
$ws2.Rows.Item(3).Delete()

# ---------------------------------------------------------------
# Sheet "Btts": 9 data rows -> 8 data rows (rows 2-10 -> rows 2-9)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Btts")

$ws3.Range("A2").Value = "19-12-2024 20:00"
$ws3.Range("B2").Value = "WORLD"
$ws3.Range("C2").Value = "UEFA EUROPA CONFERENCE LEAGUE"
$ws3.Range("D2").Value = "Larne - Gent"
$ws3.Range("E2").Value = 76
$ws3.Range("F2").Value = 2

$ws3.Range("A3").Value = "19-12-2024 20:00"
$ws3.Range("B3").Value = "WORLD"
$ws3.Range("C3").Value = "UEFA EUROPA CONFERENCE LEAGUE"
$ws3.Range("D3").Value = "Celje - The New Saints"
$ws3.Range("E3").Value = 76
$ws3.Range("F3").Value = 1.73

$ws3.Range("A4").Value = "19-12-2024 14:00"
$ws3.Range("B4").Value = "ALGERIA"
$ws3.Range("C4").Value = "U21 LEAGUE 1"
$ws3.Range("D4").Value = "CS Constantine U21 - NC Magra U21"
$ws3.Range("E4").Value = 78.3
$ws3.Range("F4").Value = 2.2

$ws3.Range("A5").Value = "19-12-2024 14:00"
$ws3.Range("B5").Value = "ALGERIA"
$ws3.Range("C5").Value = "U21 LEAGUE 1"
$ws3.Range("D5").Value = "ES Sétif U21 - USM Alger U21"
$ws3.Range("E5").Value = 76.7
$ws3.Range("F5").Value = 1.8

$ws3.Range("A6").Value = "19-12-2024 12:30"
$ws3.Range("B6").Value = "ALGERIA"
$ws3.Range("C6").Value = "U21 LEAGUE 1"
$ws3.Range("D6").Value = "MC Alger U21 - ASO Chlef U21"
$ws3.Range("E6").Value = 75.8
$ws3.Range("F6").Value = 1.8

$ws3.Range("A7").Value = "19-12-2024 20:15"
$ws3.Range("B7").Value = "SPAIN"
$ws3.Range("C7").Value = "SEGUNDA DIVISIÓN"
$ws3.Range("D7").Value = "Huesca - Tenerife"
$ws3.Range("E7").Value = 76.7
$ws3.Range("F7").Value = 2.37

$ws3.Range("A8").Value = "20-12-2024 19:45"
$ws3.Range("B8").Value = "ENGLAND"
$ws3.Range("C8").Value = "LEAGUE ONE"
$ws3.Range("D8").Value = "Cambridge United - Huddersfield"
$ws3.Range("E8").Value = 76
$ws3.Range("F8").Value = 1.8

$ws3.Range("A9").Value = "20-12-2024 19:30"
$ws3.Range("B9").Value = "ITALY"
$ws3.Range("C9").Value = "SERIE C - GIRONE B"
$ws3.Range("D9").Value = "Pontedera - Legnago Salus"
$ws3.Range("E9").Value = 88
$ws3.Range("F9").Value = 1.85

$ws3.Rows.Item(10).Delete()

# ---------------------------------------------------------------
# Sheet "Over_Under": 4 data rows -> 7 data rows (rows 2-5 -> rows 2-8)
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Over_Under")

$ws4.Range("A2").Value = "19-12-2024 20:00"
$ws4.Range("B2").Value = "WORLD"
$ws4.Range("C2").Value = "UEFA EUROPA CONFERENCE LEAGUE"
$ws4.Range("D2").Value = "1. FC Heidenheim - FC ST. Gallen"
$ws4.Range("E2").Value = 70
$ws4.Range("F2").Value = 1.53
$ws4.Range("G2").Value = 60
$ws4.Range("H2").Value = 2.25

$ws4.Range("A3").Value = "19-12-2024 20:00"
$ws4.Range("B3").Value = "WORLD"
$ws4.Range("C3").Value = "UEFA EUROPA CONFERENCE LEAGUE"
$ws4.Range("D3").Value = "Djurgardens IF - Legia Warszawa"
$ws4.Range("E3").Value = 100
$ws4.Range("F3").Value = 1.85
$ws4.Range("G3").Value = 40
$ws4.Range("H3").Value = 3

$ws4.Range("A4").Value = "19-12-2024 20:00"
$ws4.Range("B4").Value = "WORLD"
$ws4.Range("C4").Value = "UEFA EUROPA CONFERENCE LEAGUE"
$ws4.Range("D4").Value = "TSC Backa Topola - FC Noah"
$ws4.Range("E4").Value = 86.7
$ws4.Range("F4").Value = 1.65
$ws4.Range("G4").Value = 60
$ws4.Range("H4").Value = 2.6

$ws4.Range("A5").Value = "20-12-2024 19:45"
$ws4.Range("B5").Value = "ITALY"
$ws4.Range("C5").Value = "SERIE A"
$ws4.Range("D5").Value = "Verona - AC Milan"
$ws4.Range("E5").Value = 85
$ws4.Range("F5").Value = 1.73
$ws4.Range("G5").Value = 55
$ws4.Range("H5").Value = 2.75

$ws4.Range("A6").Value = "20-12-2024 19:00"
$ws4.Range("B6").Value = "NETHERLANDS"
$ws4.Range("C6").Value = "EERSTE DIVISIE"
$ws4.Range("D6").Value = "Dordrecht - FC Eindhoven"
$ws4.Range("E6").Value = 80
$ws4.Range("F6").Value = 1.48
$ws4.Range("G6").Value = 80
$ws4.Range("H6").Value = 2.2

$ws4.Range("A7").Value = "20-12-2024 19:00"
$ws4.Range("B7").Value = "BELGIUM"
$ws4.Range("C7").Value = "CHALLENGER PRO LEAGUE"
$ws4.Range("D7").Value = "Club Brugge II - Zulte Waregem"
$ws4.Range("E7").Value = 85
$ws4.Range("F7").Value = 1.7
$ws4.Range("G7").Value = 65
$ws4.Range("H7").Value = 2.75

$ws4.Range("A8").Value = "20-12-2024 13:00"
$ws4.Range("B8").Value = "WORLD"
$ws4.Range("C8").Value = "AFF CHAMPIONSHIP"
$ws4.Range("D8").Value = "Malaysia - Singapore"
$ws4.Range("E8").Value = 75
$ws4.Range("F8").Value = 1.65
$ws4.Range("G8").Value = 60
$ws4.Range("H8").Value = 2.63
